# Apply cryptos list update (coin rankings, prices, 1h volume deltas)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.299.92"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.932.10"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'0.7507"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.73%  "

$ws.Range("D6").Value = "'242.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.46%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'27.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("D9").Value = "'0.3173"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("D10").Value = "'0.07094"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("D12").Value = "'0.7787"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("D13").Value = "1.915.05"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("D14").Value = "'5.397"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").Value = "'93.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").Value = "'14.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("D17").Value = "30.297.33"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").Value = "'6.041"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.70%  "

$ws.Range("D19").Value = "'251.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.97%  "

$ws.Range("D20").Value = "'0.000007944"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.161.41"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  -2.17%  "

$ws.Range("D25").Value = "'9.571"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").Value = "'165.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("D28").Value = "'0.1297"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.13%  "

$ws.Range("D29").Value = "'2.186"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.76%  "

$ws.Range("D30").Value = "'1.368"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.03%  "

$ws.Range("D31").Value = "'1.557"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.96%  "

$ws.Range("D32").Value = "'4.415"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("D33").Value = "'4.145"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("D34").Value = "'0.05235"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").Value = "'1.319"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.99%  "

$ws.Range("D37").Value = "'2.777"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.48%  "

$ws.Range("D38").Value = "'0.01954"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").Value = "'6.514"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.43%  "

$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("D42").Value = "'0.4531"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").Value = "'1.983"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D44").Value = "'0.8412"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.689"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.37%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.978"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.11%  "

$ws.Range("D48").Value = "'101.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.14%  "

$ws.Range("D49").Value = "'38.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.63%  "

$ws.Range("D50").Value = "'0.1233"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.61%  "

$ws.Range("D51").Value = "'958.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.61%  "

